$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B (Coin), C (Link), D (Price), E (Volume(1h)) are treated as text
# so that numeric-looking strings (e.g. "247.04") are preserved exactly as text
# and not auto-converted into numbers by Excel.
$ws.Columns.Item(2).NumberFormat = "@"
$ws.Columns.Item(3).NumberFormat = "@"
$ws.Columns.Item(4).NumberFormat = "@"
$ws.Columns.Item(5).NumberFormat = "@"

$ws.Range("D2").Value = "42.746.85"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "2.258.76"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "247.04"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "77.31"
$ws.Range("E7").Value = "  +2.70%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.628"
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("D10").Value = "45.58"
$ws.Range("E10").Value = "  +12.71%  "
$ws.Range("D11").Value = "0.0956"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "7.36"
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "14.81"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "0.866"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "2.259.51"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "42.485.35"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("E18").Value = "  +4.11%  "
$ws.Range("D19").Value = "6.23"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "72.30"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "10.21"
$ws.Range("E21").Value = "  +40.43%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Value = "2.28"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "233.18"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "11.91"
$ws.Range("E24").Value = "  +5.15%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "3.63"
$ws.Range("E26").Value = "  -2.43%  "
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").Value = "2.23"
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("D29").Value = "167.36"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("D30").Value = "20.79"
$ws.Range("E30").Value = "  +1.06%  "
$ws.Range("D31").Value = "0.0832"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "31.85"
$ws.Range("E32").Value = "  -6.82%  "
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "5.46"
$ws.Range("E34").Value = "  +11.03%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").Value = "4.70"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("E37").Value = "  +5.36%  "
$ws.Range("D38").Value = "14.59"
$ws.Range("E38").Value = "  +7.42%  "
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("D40").Value = "5.84"
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").Value = "64.46"
$ws.Range("E41").Value = "  +6.73%  "
$ws.Range("D42").Value = "0.204"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").Value = "108.85"
$ws.Range("E43").Value = "  -1.53%  "
$ws.Range("D44").Value = "8.92"
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("E45").Value = "  +3.32%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "2.42"
$ws.Range("E47").Value = "  +6.52%  "
$ws.Range("D48").Value = "1.15"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D50").Value = "4.18"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("E51").Value = "  +0.50%  "
